# Kujata_Profits market-data refresh ('chore: update Sheets via scheduled runner').
#
# The scheduled runner recomputes, per Leve row, the live marketboard columns
#   H  currentAveragePrice     I  currentAveragePriceNQ   J  currentAveragePriceHQ
#   K  LevePriceNQ              L  LevePriceHQ
#   M  LeveProfitNQ             N  LeveProfitHQ
# from fresh market prices. A handful of rows also lose their NQ- or HQ-side
# profit cell entirely once that side's listings dry up (no HQ/NQ data -> the
# corresponding LeveProfit cell is cleared rather than zeroed), so those use
# ClearContents() to match.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 111.111115
$ws.Range("I28").Value = 111.111115
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 111.111115
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 373.888885
$ws.Range("N28").ClearContents()
# Row 29: Dripping with Venom
$ws.Range("H29").Value = 1990.909
$ws.Range("I29").Value = 250
$ws.Range("J29").Value = 2377.7778
$ws.Range("K29").Value = 750
$ws.Range("L29").Value = 7133.3334
$ws.Range("M29").Value = -469
$ws.Range("N29").Value = -7695.3334
# Row 32: Automata for the People
$ws.Range("H32").Value = 1587.3
$ws.Range("I32").Value = 1400
$ws.Range("J32").Value = 1620.3529
$ws.Range("K32").Value = 1400
$ws.Range("L32").Value = 1620.3529
$ws.Range("M32").Value = -1074
$ws.Range("N32").Value = -2272.3529
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 9526216
$ws.Range("I132").Value = 10103345
$ws.Range("J132").Value = 3590
$ws.Range("K132").Value = 30310035
$ws.Range("L132").Value = 10770
$ws.Range("M132").Value = -30307505
$ws.Range("N132").Value = -15830
# Row 134: Binding Spells
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 135: For Tired Minds
$ws.Range("H135").Value = 800
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 800
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 7200
$ws.Range("N135").Value = -12270
$ws.Range("M135").ClearContents()
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1231.3939
$ws.Range("I137").Value = 1178.8636
$ws.Range("J137").Value = 1336.4546
$ws.Range("K137").Value = 3536.5908
$ws.Range("L137").Value = 4009.3638
$ws.Range("M137").Value = -986.5907999999999
$ws.Range("N137").Value = -9109.363799999999
# Row 140: Tome for Tradition
$ws.Range("H140").Value = 30675
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 30675
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 30675
$ws.Range("N140").Value = -41035

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 8128.357
$ws.Range("I2").Value = 1011
$ws.Range("J2").Value = 17618.166
$ws.Range("K2").Value = 1011
$ws.Range("L2").Value = 17618.166
$ws.Range("M2").Value = -898
$ws.Range("N2").Value = -17844.166
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2313.5
$ws.Range("I61").Value = 1907.3334
$ws.Range("J61").Value = 2487.5715
$ws.Range("K61").Value = 1907.3334
$ws.Range("L61").Value = 2487.5715
$ws.Range("M61").Value = -1695.3334
$ws.Range("N61").Value = -2911.5715
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1047.6086
$ws.Range("I74").Value = 863.58826
$ws.Range("J74").Value = 1569
$ws.Range("K74").Value = 863.58826
$ws.Range("L74").Value = 1569
$ws.Range("M74").Value = 10.41174000000001
$ws.Range("N74").Value = -3317
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1047.6086
$ws.Range("I77").Value = 863.58826
$ws.Range("J77").Value = 1569
$ws.Range("K77").Value = 4317.9413
$ws.Range("L77").Value = 7845
$ws.Range("M77").Value = 50.05869999999959
$ws.Range("N77").Value = -16581
# Row 116: No Scope
$ws.Range("H116").Value = 8128.357
$ws.Range("I116").Value = 1011
$ws.Range("J116").Value = 17618.166
$ws.Range("K116").Value = 1011
$ws.Range("L116").Value = 17618.166
$ws.Range("M116").Value = 1283
$ws.Range("N116").Value = -22206.166
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2313.5
$ws.Range("I136").Value = 1907.3334
$ws.Range("J136").Value = 2487.5715
$ws.Range("K136").Value = 5722.0002
$ws.Range("L136").Value = 7462.7145
$ws.Range("M136").Value = -3172.0002
$ws.Range("N136").Value = -12562.7145
# Row 140: A Hand for a Deckhand
$ws.Range("H140").Value = 73180
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 73180
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 73180
$ws.Range("N140").Value = -83540

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 8128.357
$ws.Range("I3").Value = 1011
$ws.Range("J3").Value = 17618.166
$ws.Range("K3").Value = 1011
$ws.Range("L3").Value = 17618.166
$ws.Range("M3").Value = -897
$ws.Range("N3").Value = -17846.166
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 2328.7646
$ws.Range("I20").Value = 2180
$ws.Range("J20").Value = 2541.2856
$ws.Range("K20").Value = 2180
$ws.Range("L20").Value = 2541.2856
$ws.Range("M20").Value = -1933
$ws.Range("N20").Value = -3035.2856
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 31251420
$ws.Range("I99").Value = 38462876
$ws.Range("J99").Value = 1766.3334
$ws.Range("K99").Value = 38462876
$ws.Range("L99").Value = 1766.3334
$ws.Range("M99").Value = -38461378
$ws.Range("N99").Value = -4762.3334
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 142859040
$ws.Range("I105").Value = 250001740
$ws.Range("J105").Value = 2100
$ws.Range("K105").Value = 250001740
$ws.Range("L105").Value = 2100
$ws.Range("M105").Value = -249999993
$ws.Range("N105").Value = -5594
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 10892.357
$ws.Range("I134").Value = 1242
$ws.Range("J134").Value = 20542.715
$ws.Range("K134").Value = 3726
$ws.Range("L134").Value = 61628.145
$ws.Range("M134").Value = -1191
$ws.Range("N134").Value = -66698.145

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1065.7428
$ws.Range("I58").Value = 934.0417
$ws.Range("J58").Value = 1353.091
$ws.Range("K58").Value = 934.0417
$ws.Range("L58").Value = 1353.091
$ws.Range("M58").Value = -731.0417
$ws.Range("N58").Value = -1759.091
# Row 60: Bowing to Greater Power
$ws.Range("H60").Value = 9238.799999999999
$ws.Range("I60").Value = 2174.1428
$ws.Range("J60").Value = 11986.167
$ws.Range("K60").Value = 2174.1428
$ws.Range("L60").Value = 11986.167
$ws.Range("M60").Value = -1663.1428
$ws.Range("N60").Value = -13008.167
# Row 107: Built to Last
$ws.Range("H107").Value = 569.96295
$ws.Range("I107").Value = 450.11765
$ws.Range("J107").Value = 773.7
$ws.Range("K107").Value = 450.11765
$ws.Range("L107").Value = 773.7
$ws.Range("M107").Value = 1469.88235
$ws.Range("N107").Value = -4613.7
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 27780032
$ws.Range("I134").Value = 41668960
$ws.Range("J134").Value = 2175
$ws.Range("K134").Value = 125006880
$ws.Range("L134").Value = 6525
$ws.Range("M134").Value = -125004345
$ws.Range("N134").Value = -11595
# Row 136: Turali Quality
$ws.Range("H136").Value = 1065.7428
$ws.Range("I136").Value = 934.0417
$ws.Range("J136").Value = 1353.091
$ws.Range("K136").Value = 2802.1251
$ws.Range("L136").Value = 4059.273
$ws.Range("M136").Value = -252.1251000000002
$ws.Range("N136").Value = -9159.272999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 17: Chew the Fat
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
# Row 34: Fever Pitch
$ws.Range("H34").Value = 1442.7059
$ws.Range("I34").Value = 421.5
$ws.Range("J34").Value = 1999.7273
$ws.Range("K34").Value = 1264.5
$ws.Range("L34").Value = 5999.1819
$ws.Range("M34").Value = -1180.5
$ws.Range("N34").Value = -6167.1819
# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 1615.1538
$ws.Range("I39").Value = 403
$ws.Range("J39").Value = 1716.1666
$ws.Range("K39").Value = 1209
$ws.Range("L39").Value = 5148.4998
$ws.Range("M39").Value = -915
$ws.Range("N39").Value = -5736.4998
# Row 55: Pagan Pastries
$ws.Range("H55").Value = 3499.6667
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3499.6667
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10499.0001
$ws.Range("N55").Value = -10853.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 22506628
$ws.Range("I70").Value = 20839178
$ws.Range("J70").Value = 25007802
$ws.Range("K70").Value = 20839178
$ws.Range("L70").Value = 25007802
$ws.Range("M70").Value = -20838908
$ws.Range("N70").Value = -25008342
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 22506628
$ws.Range("I73").Value = 20839178
$ws.Range("J73").Value = 25007802
$ws.Range("K73").Value = 20839178
$ws.Range("L73").Value = 25007802
$ws.Range("M73").Value = -20838242
$ws.Range("N73").Value = -25009674
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 3756.4546
$ws.Range("I80").Value = 2316.4285
$ws.Range("J80").Value = 6276.5
$ws.Range("K80").Value = 2316.4285
$ws.Range("L80").Value = 6276.5
$ws.Range("M80").Value = -1318.4285
$ws.Range("N80").Value = -8272.5
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 3756.4546
$ws.Range("I83").Value = 2316.4285
$ws.Range("J83").Value = 6276.5
$ws.Range("K83").Value = 11582.1425
$ws.Range("L83").Value = 31382.5
$ws.Range("M83").Value = -6590.1425
$ws.Range("N83").Value = -41366.5
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 216485.72
$ws.Range("I122").Value = 2780
$ws.Range("J122").Value = 750750
$ws.Range("K122").Value = 8340
$ws.Range("L122").Value = 2252250
$ws.Range("M122").Value = -5890
$ws.Range("N122").Value = -2257150
# Row 132: On Board for Lar
$ws.Range("H132").Value = 3301.125
$ws.Range("I132").Value = 3524.6
$ws.Range("J132").Value = 3199.5454
$ws.Range("K132").Value = 10573.8
$ws.Range("L132").Value = 9598.636200000001
$ws.Range("M132").Value = -8043.799999999999
$ws.Range("N132").Value = -14658.6362

$ws = $wb.Worksheets.Item("LTW")
# Row 42: Slave to Fashion
$ws.Range("H42").Value = 5850
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 5850
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 5850
$ws.Range("N42").Value = -6976
# Row 49: First They Came for the Heretics
$ws.Range("H49").Value = 5850
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 5850
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 5850
$ws.Range("N49").Value = -6144
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 980
$ws.Range("I61").Value = 725
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 725
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -523
$ws.Range("N61").Value = -2404
# Row 113: Peace in Rest
$ws.Range("H113").Value = 980
$ws.Range("I113").Value = 725
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 725
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1445
$ws.Range("N113").Value = -6340
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 13164.777
$ws.Range("I136").Value = 34994.668
$ws.Range("J136").Value = 2249.8333
$ws.Range("K136").Value = 104984.004
$ws.Range("L136").Value = 6749.499899999999
$ws.Range("M136").Value = -102434.004
$ws.Range("N136").Value = -11849.4999

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 5427.143
$ws.Range("I132").Value = 4451.5
$ws.Range("J132").Value = 5817.4
$ws.Range("K132").Value = 13354.5
$ws.Range("L132").Value = 17452.2
$ws.Range("M132").Value = -10824.5
$ws.Range("N132").Value = -22512.2
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 960.75
$ws.Range("I136").Value = 960.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2882.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -332.25
$ws.Range("N136").ClearContents()
# Row 138: Halfgloves, Full Effort
$ws.Range("H138").Value = 34963
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 34963
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 34963
$ws.Range("N138").Value = -45243

Write-Output "Kujata_Profits: refreshed market data for 43 Leve rows across 8 sheets."
